# Fruta / hortaliza, semanal
# Insert two new weekly data rows at the top of the data block (rows 14-15),
# pushing the existing rows 14-36 down to 16-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 14 (shifts old rows 14..36 down to 16..38)
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# New row 14: Packham's Triumph, Primera, 2022-08-18
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = "2022-08-18"
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100104
$ws.Range("H14").Value = "Frutos de pepita"
$ws.Range("I14").Value = 100104005
$ws.Range("J14").Value = "Pera"
$ws.Range("K14").Value = "Packham's Triumph"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 19000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19500
$ws.Range("Q14").Value = "`$/caja 18 kilos granel"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 1083
$ws.Range("T14").Value = 18

# New row 15: Winter Nelis, Primera, 2022-08-18
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C15").Value = "Arica y Parinacota"
$ws.Range("D15").Value = "2022-08-18"
$ws.Range("E15").Value = 15
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100104
$ws.Range("H15").Value = "Frutos de pepita"
$ws.Range("I15").Value = 100104005
$ws.Range("J15").Value = "Pera"
$ws.Range("K15").Value = "Winter Nelis"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 270
$ws.Range("N15").Value = 19000
$ws.Range("O15").Value = 20000
$ws.Range("P15").Value = 19500
$ws.Range("Q15").Value = "`$/caja 18 kilos granel"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1083
$ws.Range("T15").Value = 18
